# Daily auto-push update: insert one new 3-hourly reading row for
# 2026/01/16 ahead of the existing 2026/01/16 rows, pushing every
# subsequent row (old rows 662-703) down by one (new rows 663-704).
#
# Row 661 (the last existing 2026/01/16 row: 2026/01/16, 金, 14, 201)
# already has exactly the date/weekday/ranking values the new row
# needs (2026/01/16, 金, *, 201) and carries no special formatting, so
# copying it and inserting the copy is the cleanest way to create the
# new row without Excel's autodetect turning the "2026/01/16" text
# into a real date serial (which happens if we type the literal string
# straight into a freshly inserted, default-formatted cell).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(661).Copy()
$ws.Rows.Item(662).Insert()

# Only the "time" column differs from the copied row 661 template.
$ws.Cells.Item(662, 3).Value = 17
